$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Booking #15 arrived late and paid downpayments for two more rooms (agencies).
# Mirror the existing row layout/format (copy row 16's formatting, incl. the
# date-formatted column C) down into the two new rows, then set the values.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 44905

$ws.Range("A16:C16").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Range("A18").Value = 15
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 44905

# Clear the marching-ants clipboard highlight and restore the cursor to the
# author's final selection.
$excel.CutCopyMode = 0
$ws.Range("E15").Select()
